$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The sheet currently ends with a note row at row 96 (column B).
# We need to insert a new data row *before* that note, pushing it to row 97,
# and fill the new row 96 with the next day's figures.
$ws.Rows.Item(96).Insert()

$ws.Cells.Item(96, 1).Value = 43951
$ws.Cells.Item(96, 2).Value = 792
$ws.Cells.Item(96, 3).Value = 32302
$ws.Cells.Item(96, 4).Value = 185
$ws.Cells.Item(96, 5).Value = 6849

foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$99"
    }
}

$ws.Range("D97").Select() | Out-Null
